$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in previously-blank ISO_3166-1 / Continente columns for a handful of
# territories/countries that had been left empty.
$ws.Range("D18").Value = "BMU"
$ws.Range("E18").Value = "América"

$ws.Range("D67").Value = "HKG"
$ws.Range("E67").Value = "Asia"

$ws.Range("D75").Value = "ISR"
$ws.Range("E75").Value = "Asia"

$ws.Range("D116").Value = "PRI"
$ws.Range("E116").Value = "América"

$ws.Range("D145").Value = "TWN"
$ws.Range("E145").Value = "Asia"

# Row 147 was a duplicate "Tajikistan (TJK)" row (identical to row 146).
# Remove it so every following row shifts up by one, eliminating the
# duplicate and realigning the rest of the table (including the duplicate
# "Turkey (TUR)" rows at 154/155, which become distinct Europe/Asia rows
# once shifted).
$ws.Rows(147).Delete()
